# "Add configuration for m/m/1"
#
# Inserts a new worksheet ("Arkusz4") right after "Arkusz1" that holds the
# raw write-count samples plus their LOG10() transform, then rewires
# "Arkusz1" to work off those log10-transformed numbers (H3:J6), adding the
# extra (x-mean)^2 columns (U:W) / grand mean (J8) / total (X7) needed for
# the new analysis.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Arkusz1")

# ------------------------------------------------------------------
# 1) Insert the new "Arkusz4" sheet right after "Arkusz1"
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws1)
$ws4.Name = "Arkusz4"

# Raw replication data (same numbers that used to live in Arkusz1!H3:J6)
$ws4.Range("A1").Value = 12833
$ws4.Range("B1").Value = 12684
$ws4.Range("C1").Value = 12974

$ws4.Range("A2").Value = 12732
$ws4.Range("B2").Value = 12038
$ws4.Range("C2").Value = 12121

$ws4.Range("A3").Value = 12140
$ws4.Range("B3").Value = 12667
$ws4.Range("C3").Value = 12533

$ws4.Range("A4").Value = 10582
$ws4.Range("B4").Value = 11077
$ws4.Range("C4").Value = 10188

# LOG10() transform of the raw samples
$ws4.Range("E1").Formula = "=LOG10(A1)"
$ws4.Range("F1:G1").Formula = "=LOG10(B1)"

$ws4.Range("E2:E4").Formula = "=LOG10(A2)"
$ws4.Range("F2:F4").Formula = "=LOG10(B2)"
$ws4.Range("G2:G4").Formula = "=LOG10(C2)"

# ------------------------------------------------------------------
# 2) Arkusz1: replace the raw write counts with their log10 values
# ------------------------------------------------------------------
$ws1.Range("H3").Value = $ws4.Range("E1").Value()
$ws1.Range("I3").Value = $ws4.Range("F1").Value()
$ws1.Range("J3").Value = $ws4.Range("G1").Value()

$ws1.Range("H4").Value = $ws4.Range("E2").Value()
$ws1.Range("I4").Value = $ws4.Range("F2").Value()
$ws1.Range("J4").Value = $ws4.Range("G2").Value()

$ws1.Range("H5").Value = $ws4.Range("E3").Value()
$ws1.Range("I5").Value = $ws4.Range("F3").Value()
$ws1.Range("J5").Value = $ws4.Range("G3").Value()

$ws1.Range("H6").Value = $ws4.Range("E4").Value()
$ws1.Range("I6").Value = $ws4.Range("F4").Value()
$ws1.Range("J6").Value = $ws4.Range("G4").Value()

# ------------------------------------------------------------------
# 3) Arkusz1: grand mean + squared-deviation-from-grand-mean columns
# ------------------------------------------------------------------
$ws1.Range("J8").Formula = "=AVERAGE(H3:J6)"

$ws1.Range("U3:U6").Formula = '=(H3-$J$8)^2'
$ws1.Range("V3:W3").Formula = '=(I3-$J$8)^2'
$ws1.Range("V4:W6").Formula = '=(I4-$J$8)^2'

$ws1.Range("X7").Formula = "=SUM(U3:W6)"

# ------------------------------------------------------------------
# 4) Selections / active sheet / tab order bookkeeping
# ------------------------------------------------------------------
# Arkusz2 (the LaTeX-export sheet) used to be the active tab with K24
# selected; it now just sits quietly on K6.
$ws2 = $wb.Worksheets.Item("Arkusz2")
$ws2.Range("K6").Select()

# Arkusz4's interesting range is its log10 table.
$ws4.Range("E1:G4").Select()

# Arkusz1 becomes the active tab, focused on the newly-populated H3:J6
# sample block.
$ws1.Activate()
$ws1.Range("H3:J6").Select()

# ------------------------------------------------------------------
# 5) Window geometry (cosmetic)
# ------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 460
$win.Width = 28800
$win.Height = 16000
